# "Generate Report for Handback" - fill in the handback results for the
# 5c90f014-e081-45c8-b275-899243d01f20 row (row 7) on both the zh-cn and
# de-de sheets: a link to the target file, the handback xlf file name, the
# handback datetime, and an error message (version mismatch).

$wb = $excel.ActiveWorkbook

$errorDetail = 'The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/15585eb81f8ff3b8fec5752882a01ca8d4895044/e2e/5c90f014-e081-45c8-b275-899243d01f20.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/f9d5d7578f8e87048658eebe9a05bd680dae7f77/e2e/5c90f014-e081-45c8-b275-899243d01f20.md.'

function Fill-HandbackRow($SheetName, $TargetUrl, $HandbackFile, $HandbackDateTime) {

    $ws = $wb.Worksheets.Item($SheetName)

    # Latest Target File: link display text + real hyperlink
    $ws.Range("I7").Value = "5c90f014-e081-45c8-b275-899243d01f20.md"
    $ws.Range("I7").Style = "Hyperlink"
    $ws.Hyperlinks.Add($ws.Range("I7"), $TargetUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "5c90f014-e081-45c8-b275-899243d01f20.md") | Out-Null

    # Latest Handback File
    $ws.Range("J7").Value = $HandbackFile

    # Latest Handback DateTime
    $ws.Range("K7").Value = $HandbackDateTime

    # Error Detail
    $ws.Range("P7").Value = $errorDetail
}

Fill-HandbackRow "zh-cn" "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f9d5d7578f8e87048658eebe9a05bd680dae7f77/e2e/5c90f014-e081-45c8-b275-899243d01f20.md" "5c90f014-e081-45c8-b275-899243d01f20.b66427909297682ee7169af699c750ea24d680a2.zh-cn.xlf" "2016-08-12 15:06:26"

Fill-HandbackRow "de-de" "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/f9d5d7578f8e87048658eebe9a05bd680dae7f77/e2e/5c90f014-e081-45c8-b275-899243d01f20.md" "5c90f014-e081-45c8-b275-899243d01f20.b66427909297682ee7169af699c750ea24d680a2.de-de.xlf" "2016-08-12 15:06:35"
